$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set final cell values (rows 1-23)
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"
$ws.Range("B2").Value = "LOB1265"
$ws.Range("C2").Value = "LOB1265"
$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Epidemiologia Ambiental"
$ws.Range("C3").Value = " Epidemiologia Ambiental"
$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Environmental Epidemiology"
$ws.Range("C4").Value = "Environmental Epidemiology"
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "2"
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"
$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EA-6"
$ws.Range("C9").Value = "EA-6"
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("C10").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Train higher education professionals with knowledge in epidemiology that allows them to act in the assessment of environmental risk and public health, assessment of environmental impacts on public health, as well as to act in the management of this risk."
$ws.Range("C11").Value = "Train higher education professionals with knowledge in epidemiology that allows them to act in the assessment of environmental risk and public health, assessment of environmental impacts on public health, as well as to act in the management of this risk."
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Definition of epidemiology and its areas of activity; Types of epidemiological studies; Errors and confounding factors in epidemiological studies; Environmental epidemiology: exposure and quantification of exposure to environmental factors; Risk assessment; Impact assessment on the environment and public health; Risk management."
$ws.Range("C14").Value = "Definition of epidemiology and its areas of activity; Types of epidemiological studies; Errors and confounding factors in epidemiological studies; Environmental epidemiology: exposure and quantification of exposure to environmental factors; Risk assessment; Impact assessment on the environment and public health; Risk management."
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Introduction: Definition, field of action of epidemiology; Epidemiology and public health. Types of studies: Observational epidemiological studies (descriptive studies, ecological or correlation studies, ecological fallacy, cross-sectional studies, case and control studies, cohort studies); Experimental epidemiological studies (randomized clinical trial, field trials, community trials). Errors and confounding factors: Potential errors in epidemiological studies (random error, sample size, systematic error, selection bias, measurement bias); Confounding factors (control of confounding factors, validity, ethical issues). Environmental epidemiology: Exposure to environmental factors and quantification of exposure: Biological monitoring; Interpretation of biological data; Individual measures versus group measures; Population dose; Dose-effect relationship and dose-response relationship. Risk: Risk assessment; Health impact assessment; Risk management; Environmental and public health impact assessment."
$ws.Range("C16").Value = "Introduction: Definition, field of action of epidemiology; Epidemiology and public health. Types of studies: Observational epidemiological studies (descriptive studies, ecological or correlation studies, ecological fallacy, cross-sectional studies, case and control studies, cohort studies); Experimental epidemiological studies (randomized clinical trial, field trials, community trials). Errors and confounding factors: Potential errors in epidemiological studies (random error, sample size, systematic error, selection bias, measurement bias); Confounding factors (control of confounding factors, validity, ethical issues). Environmental epidemiology: Exposure to environmental factors and quantification of exposure: Biological monitoring; Interpretation of biological data; Individual measures versus group measures; Population dose; Dose-effect relationship and dose-response relationship. Risk: Risk assessment; Health impact assessment; Risk management; Environmental and public health impact assessment."
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("C18").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas teóricas expositivas com resolução de exercícios e discussão de casos reais de impactos ambientais e seus potenciais reflexos à saúde pública: análise de riscos, avaliação dos impactos ambientais, avaliação dos impactos à saúde pública e manejo de riscos."
$ws.Range("C19").Value = "Aulas teóricas expositivas com resolução de exercícios e discussão de casos reais de impactos ambientais e seus potenciais reflexos à saúde pública: análise de riscos, avaliação dos impactos ambientais, avaliação dos impactos à saúde pública e manejo de riscos."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "O sistema de avaliação será composto por 2 avaliações de igual peso. A Nota Final será obtida por meio da média simples dessas duas avaliações. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.(Nota final+P_recuperação)/2"
$ws.Range("C20").Value = "O sistema de avaliação será composto por 2 avaliações de igual peso. A Nota Final será obtida por meio da média simples dessas duas avaliações. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.(Nota final+P_recuperação)/2"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Estará em período de recuperação o aluno que obtiver notas entre 3,0 e 4,9. Para esses alunos a Nota Finalrec será calculada pela média simples entre a avaliação de recuperação (todo o conteúdo do semestre) e sua nota final.(Nota final+P_recuperação)/2"
$ws.Range("C21").Value = "Estará em período de recuperação o aluno que obtiver notas entre 3,0 e 4,9. Para esses alunos a Nota Finalrec será calculada pela média simples entre a avaliação de recuperação (todo o conteúdo do semestre) e sua nota final.(Nota final+P_recuperação)/2"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"

# Clear stray cells that existed in the old layout but must be empty now
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# Remove now-unused row 24 entirely (dimension shrinks to A1:C23)
$ws.Rows.Item(24).Delete()

# Fix row heights
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
